$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the short-code column (C) actually changes value; the domain
# column (B) keeps the same text on every row, so it is left untouched.
$ws.Range("C1").Value = "37w8pj6m"
$ws.Range("C2").Value = "5tkkcopb"
$ws.Range("C3").Value = "v6zhob9m"
